$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing header ("时间") and its sample value that had been
# misplaced on import - they belong in column B, next to the existing
# "姓名" column.
$ws.Range("B1").Value = "时间"
$ws.Range("B1").Font.Name = "宋体"
$ws.Range("B1").Font.Size = 11
$ws.Range("B1").WrapText = $true

$ws.Range("B2").Value = 20190901

# Size column B to fit its new content (matches the best-fit width Excel
# computed for the "时间"/date column).
$ws.Columns("B:B").ColumnWidth = 8.9

# Restore a sane selection (was left on the old, now out-of-place F14).
$ws.Range("B4").Select()

# Give the sheet explicit page setup (paper size / orientation) as it now
# has real content worth printing.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
